$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and C (move B to after C, i.e. swap positions 2 and 3)
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(4).Insert() | Out-Null

# Swap columns D and E (move D to after E, i.e. swap positions 4 and 5)
$ws.Columns.Item(4).Cut() | Out-Null
$ws.Columns.Item(6).Insert() | Out-Null

# Select entire column D (matches the author's final on-screen selection)
$ws.Columns.Item(4).Select() | Out-Null
